$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell replacements (1-indexed rows in the lone-column table)
$t.Cell(1,1).Range.Text  = "0M"
$t.Cell(2,1).Range.Text  = "0M"
$t.Cell(3,1).Range.Text  = "0M"
$t.Cell(4,1).Range.Text  = "6059"
$t.Cell(6,1).Range.Text  = "0.02882"
$t.Cell(7,1).Range.Text  = "0.00820"
$t.Cell(8,1).Range.Text  = "0.00015"
$t.Cell(9,1).Range.Text  = "0.02882"
$t.Cell(10,1).Range.Text = "0.02882"
$t.Cell(11,1).Range.Text = "0.02882"
$t.Cell(12,1).Range.Text = "1.21759"

# The last three rows previously held a whole tab-separated stats line
# (count, mean, median, ... , percentage) crammed into a single run; they
# now collapse down to just the single summary value that used to live in
# row 1/2/3.
$t.Cell(44,1).Range.Text = "99.97"
$t.Cell(45,1).Range.Text = "1.22"
$t.Cell(46,1).Range.Text = "3849"
